# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the e78d4cf2-... row (row 4) on both the zh-cn and de-de sheets,
# as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-26 05:31:01"
$wsZhCn.Range("G4").Value = "2016-01-26 05:31:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-26 05:31:12"
$wsDeDe.Range("G4").Value = "2016-01-26 05:32:03"
